# Fixed update to excel issue
#
# - Rename the "Requested quantity" header on the existing sheets to more
#   descriptive names (Weekly_PO_Qty / Monthly_PO_Qty).
# - Add a new "PO Forecast" sheet containing the forecast output
#   (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

# --- Update header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets (0.75in/0.75in/1in/1in, 0.5in header/footer)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows: ds, PO_Forecast, yhat_lower, yhat_upper ---
$wsForecast.Range("A2").Value = 45018.99999999999; $wsForecast.Range("B2").Value = 269; $wsForecast.Range("C2").Value = -714.2234298111636; $wsForecast.Range("D2").Value = 1260.204545759068
$wsForecast.Range("A3").Value = 45025.99999999999; $wsForecast.Range("B3").Value = 279; $wsForecast.Range("C3").Value = -672.2827752808158; $wsForecast.Range("D3").Value = 1270.86607654051
$wsForecast.Range("A4").Value = 45032.99999999999; $wsForecast.Range("B4").Value = 288; $wsForecast.Range("C4").Value = -683.5878980682875; $wsForecast.Range("D4").Value = 1219.703128046521
$wsForecast.Range("A5").Value = 45039.99999999999; $wsForecast.Range("B5").Value = 297; $wsForecast.Range("C5").Value = -616.8774102187695; $wsForecast.Range("D5").Value = 1242.674475181399
$wsForecast.Range("A6").Value = 45046.99999999999; $wsForecast.Range("B6").Value = 306; $wsForecast.Range("C6").Value = -677.7600592432736; $wsForecast.Range("D6").Value = 1196.98558147814
$wsForecast.Range("A7").Value = 45053.99999999999; $wsForecast.Range("B7").Value = 316; $wsForecast.Range("C7").Value = -583.0451199805439; $wsForecast.Range("D7").Value = 1266.71500959079
$wsForecast.Range("A8").Value = 45060.99999999999; $wsForecast.Range("B8").Value = 325; $wsForecast.Range("C8").Value = -693.0575956890119; $wsForecast.Range("D8").Value = 1220.765066900166
$wsForecast.Range("A9").Value = 45067.99999999999; $wsForecast.Range("B9").Value = 334; $wsForecast.Range("C9").Value = -618.6886710539609; $wsForecast.Range("D9").Value = 1363.033100096997
$wsForecast.Range("A10").Value = 45074.99999999999; $wsForecast.Range("B10").Value = 343; $wsForecast.Range("C10").Value = -570.1242031776339; $wsForecast.Range("D10").Value = 1303.536053785661
$wsForecast.Range("A11").Value = 45088.99999999999; $wsForecast.Range("B11").Value = 362; $wsForecast.Range("C11").Value = -613.7829214516798; $wsForecast.Range("D11").Value = 1349.342796678671
$wsForecast.Range("A12").Value = 45102.99999999999; $wsForecast.Range("B12").Value = 380; $wsForecast.Range("C12").Value = -575.7534713826835; $wsForecast.Range("D12").Value = 1397.136981479546
$wsForecast.Range("A13").Value = 45109.99999999999; $wsForecast.Range("B13").Value = 389; $wsForecast.Range("C13").Value = -584.6774220690438; $wsForecast.Range("D13").Value = 1375.041212081185
$wsForecast.Range("A14").Value = 45116.99999999999; $wsForecast.Range("B14").Value = 399; $wsForecast.Range("C14").Value = -534.8802176145153; $wsForecast.Range("D14").Value = 1373.34252862427
$wsForecast.Range("A15").Value = 45144.99999999999; $wsForecast.Range("B15").Value = 435; $wsForecast.Range("C15").Value = -489.1484321558914; $wsForecast.Range("D15").Value = 1422.593228440426
$wsForecast.Range("A16").Value = 45151.99999999999; $wsForecast.Range("B16").Value = 445; $wsForecast.Range("C16").Value = -559.5290786739164; $wsForecast.Range("D16").Value = 1337.096440960379
$wsForecast.Range("A17").Value = 45165.99999999999; $wsForecast.Range("B17").Value = 463; $wsForecast.Range("C17").Value = -483.9122693424816; $wsForecast.Range("D17").Value = 1408.382683154335
$wsForecast.Range("A18").Value = 45172.99999999999; $wsForecast.Range("B18").Value = 472; $wsForecast.Range("C18").Value = -535.076424382009; $wsForecast.Range("D18").Value = 1405.694828394207
$wsForecast.Range("A19").Value = 45179.99999999999; $wsForecast.Range("B19").Value = 482; $wsForecast.Range("C19").Value = -410.3134886090937; $wsForecast.Range("D19").Value = 1403.028096229173
$wsForecast.Range("A20").Value = 45186.99999999999; $wsForecast.Range("B20").Value = 491; $wsForecast.Range("C20").Value = -502.1510812138894; $wsForecast.Range("D20").Value = 1469.826674934353
$wsForecast.Range("A21").Value = 45193.99999999999; $wsForecast.Range("B21").Value = 500; $wsForecast.Range("C21").Value = -485.6359335828119; $wsForecast.Range("D21").Value = 1472.92340596937
$wsForecast.Range("A22").Value = 45200.99999999999; $wsForecast.Range("B22").Value = 509; $wsForecast.Range("C22").Value = -463.9489851508926; $wsForecast.Range("D22").Value = 1433.730653123324
$wsForecast.Range("A23").Value = 45207.99999999999; $wsForecast.Range("B23").Value = 518; $wsForecast.Range("C23").Value = -392.551052713534; $wsForecast.Range("D23").Value = 1452.791149230966
$wsForecast.Range("A24").Value = 45214.99999999999; $wsForecast.Range("B24").Value = 528; $wsForecast.Range("C24").Value = -426.7573216657827; $wsForecast.Range("D24").Value = 1409.128298497946
$wsForecast.Range("A25").Value = 45221.99999999999; $wsForecast.Range("B25").Value = 537; $wsForecast.Range("C25").Value = -465.055227568831; $wsForecast.Range("D25").Value = 1516.393271457269
$wsForecast.Range("A26").Value = 45249.99999999999; $wsForecast.Range("B26").Value = 574; $wsForecast.Range("C26").Value = -342.8508270749249; $wsForecast.Range("D26").Value = 1572.843577569103
$wsForecast.Range("A27").Value = 45270.99999999999; $wsForecast.Range("B27").Value = 602; $wsForecast.Range("C27").Value = -409.2777841054047; $wsForecast.Range("D27").Value = 1535.645558920276
$wsForecast.Range("A28").Value = 45326.99999999999; $wsForecast.Range("B28").Value = 675; $wsForecast.Range("C28").Value = -272.5053783645934; $wsForecast.Range("D28").Value = 1574.932782835165
$wsForecast.Range("A29").Value = 45333.99999999999; $wsForecast.Range("B29").Value = 685; $wsForecast.Range("C29").Value = -219.4326647286063; $wsForecast.Range("D29").Value = 1746.190556170372
$wsForecast.Range("A30").Value = 45340.99999999999; $wsForecast.Range("B30").Value = 694; $wsForecast.Range("C30").Value = -284.7143924501513; $wsForecast.Range("D30").Value = 1575.088345354707
$wsForecast.Range("A31").Value = 45354.99999999999; $wsForecast.Range("B31").Value = 712; $wsForecast.Range("C31").Value = -195.8934722475439; $wsForecast.Range("D31").Value = 1681.314319347736
$wsForecast.Range("A32").Value = 45361.99999999999; $wsForecast.Range("B32").Value = 721; $wsForecast.Range("C32").Value = -283.1921889546653; $wsForecast.Range("D32").Value = 1669.896041042608
$wsForecast.Range("A33").Value = 45368.99999999999; $wsForecast.Range("B33").Value = 731; $wsForecast.Range("C33").Value = -233.4884949080786; $wsForecast.Range("D33").Value = 1670.401460183565
$wsForecast.Range("A34").Value = 45375.99999999999; $wsForecast.Range("B34").Value = 740; $wsForecast.Range("C34").Value = -196.0535318837299; $wsForecast.Range("D34").Value = 1736.652700422401
$wsForecast.Range("A35").Value = 45382.99999999999; $wsForecast.Range("B35").Value = 749; $wsForecast.Range("C35").Value = -244.6577871440195; $wsForecast.Range("D35").Value = 1702.503017309441
$wsForecast.Range("A36").Value = 45389.99999999999; $wsForecast.Range("B36").Value = 758; $wsForecast.Range("C36").Value = -158.5497232502752; $wsForecast.Range("D36").Value = 1718.504771260644
$wsForecast.Range("A37").Value = 45396.99999999999; $wsForecast.Range("B37").Value = 768; $wsForecast.Range("C37").Value = -179.2254509729821; $wsForecast.Range("D37").Value = 1693.213738648561
$wsForecast.Range("A38").Value = 45403.99999999999; $wsForecast.Range("B38").Value = 777; $wsForecast.Range("C38").Value = -96.63459680471065; $wsForecast.Range("D38").Value = 1734.809705533976
$wsForecast.Range("A39").Value = 45410.99999999999; $wsForecast.Range("B39").Value = 786; $wsForecast.Range("C39").Value = -109.7065249618074; $wsForecast.Range("D39").Value = 1761.398996001593
$wsForecast.Range("A40").Value = 45417.99999999999; $wsForecast.Range("B40").Value = 795; $wsForecast.Range("C40").Value = -137.0305099292912; $wsForecast.Range("D40").Value = 1817.503728609302
$wsForecast.Range("A41").Value = 45452.99999999999; $wsForecast.Range("B41").Value = 841; $wsForecast.Range("C41").Value = -143.7480139294747; $wsForecast.Range("D41").Value = 1770.909628661641
$wsForecast.Range("A42").Value = 45459.99999999999; $wsForecast.Range("B42").Value = 851; $wsForecast.Range("C42").Value = -103.0535370336918; $wsForecast.Range("D42").Value = 1784.499694312679
$wsForecast.Range("A43").Value = 45466.99999999999; $wsForecast.Range("B43").Value = 860; $wsForecast.Range("C43").Value = -119.9411447024792; $wsForecast.Range("D43").Value = 1834.63005946612
$wsForecast.Range("A44").Value = 45473.99999999999; $wsForecast.Range("B44").Value = 869; $wsForecast.Range("C44").Value = -127.3391982736691; $wsForecast.Range("D44").Value = 1803.504757166526
$wsForecast.Range("A45").Value = 45501.99999999999; $wsForecast.Range("B45").Value = 906; $wsForecast.Range("C45").Value = -48.64819813050064; $wsForecast.Range("D45").Value = 1883.926326595228
$wsForecast.Range("A46").Value = 45508.99999999999; $wsForecast.Range("B46").Value = 915; $wsForecast.Range("C46").Value = -26.68714095007923; $wsForecast.Range("D46").Value = 1936.38813915923
$wsForecast.Range("A47").Value = 45529.99999999999; $wsForecast.Range("B47").Value = 943; $wsForecast.Range("C47").Value = 6.153894899588735; $wsForecast.Range("D47").Value = 1893.285177741475
$wsForecast.Range("A48").Value = 45536.99999999999; $wsForecast.Range("B48").Value = 952; $wsForecast.Range("C48").Value = 29.08141101071056; $wsForecast.Range("D48").Value = 1891.782739841411
$wsForecast.Range("A49").Value = 45543.99999999999; $wsForecast.Range("B49").Value = 961; $wsForecast.Range("C49").Value = -46.69312384199959; $wsForecast.Range("D49").Value = 1924.43668139979
$wsForecast.Range("A50").Value = 45550.99999999999; $wsForecast.Range("B50").Value = 971; $wsForecast.Range("C50").Value = 2.767521734140535; $wsForecast.Range("D50").Value = 1993.514910296143
$wsForecast.Range("A51").Value = 45557.99999999999; $wsForecast.Range("B51").Value = 980; $wsForecast.Range("C51").Value = 10.01213263102387; $wsForecast.Range("D51").Value = 1906.115943747939
$wsForecast.Range("A52").Value = 45564.99999999999; $wsForecast.Range("B52").Value = 989; $wsForecast.Range("C52").Value = 49.73338632178545; $wsForecast.Range("D52").Value = 1956.131666099288
$wsForecast.Range("A53").Value = 45571.99999999999; $wsForecast.Range("B53").Value = 998; $wsForecast.Range("C53").Value = 81.84549276563546; $wsForecast.Range("D53").Value = 2025.054986862888
$wsForecast.Range("A54").Value = 45578.99999999999; $wsForecast.Range("B54").Value = 1007; $wsForecast.Range("C54").Value = 62.47751217112976; $wsForecast.Range("D54").Value = 1933.96472705401
$wsForecast.Range("A55").Value = 45585.99999999999; $wsForecast.Range("B55").Value = 1017; $wsForecast.Range("C55").Value = 33.13813203361629; $wsForecast.Range("D55").Value = 2030.688890757654
$wsForecast.Range("A56").Value = 45592.99999999999; $wsForecast.Range("B56").Value = 1026; $wsForecast.Range("C56").Value = 92.02384962717767; $wsForecast.Range("D56").Value = 1980.94388960509
$wsForecast.Range("A57").Value = 45599.99999999999; $wsForecast.Range("B57").Value = 1035; $wsForecast.Range("C57").Value = 81.96577628577091; $wsForecast.Range("D57").Value = 1951.345015456511
$wsForecast.Range("A58").Value = 45606.99999999999; $wsForecast.Range("B58").Value = 1044; $wsForecast.Range("C58").Value = 93.49485413195815; $wsForecast.Range("D58").Value = 2000.165506312068
$wsForecast.Range("A59").Value = 45613.99999999999; $wsForecast.Range("B59").Value = 1054; $wsForecast.Range("C59").Value = 91.80413084008336; $wsForecast.Range("D59").Value = 2083.468819175937
$wsForecast.Range("A60").Value = 45620.99999999999; $wsForecast.Range("B60").Value = 1063; $wsForecast.Range("C60").Value = 142.9290735039483; $wsForecast.Range("D60").Value = 1996.459465042783
$wsForecast.Range("A61").Value = 45627.99999999999; $wsForecast.Range("B61").Value = 1072; $wsForecast.Range("C61").Value = 157.7761406137423; $wsForecast.Range("D61").Value = 1984.460557607
$wsForecast.Range("A62").Value = 45634.99999999999; $wsForecast.Range("B62").Value = 1081; $wsForecast.Range("C62").Value = 111.7053804383757; $wsForecast.Range("D62").Value = 2058.827358193181
$wsForecast.Range("A63").Value = 45641.99999999999; $wsForecast.Range("B63").Value = 1091; $wsForecast.Range("C63").Value = 174.8518941149963; $wsForecast.Range("D63").Value = 2004.387515156129
$wsForecast.Range("A64").Value = 45648.99999999999; $wsForecast.Range("B64").Value = 1100; $wsForecast.Range("C64").Value = 198.1142255083013; $wsForecast.Range("D64").Value = 2078.855291821087
$wsForecast.Range("A65").Value = 45655.99999999999; $wsForecast.Range("B65").Value = 1109; $wsForecast.Range("C65").Value = 83.84836187328416; $wsForecast.Range("D65").Value = 2071.946913459396

# --- Match the formatting conventions used on the other sheets ---
# Header row: bold / centered / bordered style (same as row 1 elsewhere)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Column A: date/time number format (same as column A elsewhere)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A65").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the originally active sheet/selection
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
